$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Goal (per the diff):
#   "Please enter the number of days between 3 and 10, inclusive:"
#        -> "Please enter the number of days, between 3 and 10, inclusive:"
#   i.e. insert a comma right after "days", splitting the original single
#   run into three runs with identical run formatting:
#       1) "Please enter the number of days"
#       2) ","
#       3) " between 3 and 10, inclusive:"
#   with the (unique, auto-tracked) "_GoBack" bookmark collapsed between
#   run 2 and run 3 -- and therefore removed from its old location right
#   after the "Analyzer" run near "---=== IPC Temperature Analyzer V2.0 ===---".
# ---------------------------------------------------------------------------

# Step 1: insert the comma via Find/Replace. This rewrites the text of the
# run in place (merging to a single clean run with the comma inserted).
$rng = $d.Content
$found = $rng.Find.Execute(
    "Please enter the number of days between", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Please enter the number of days, between", 2)

# Step 2: locate the boundary right after "days" (before the comma) and
# force a run split there using a throwaway bookmark, then delete the
# bookmark -- the split survives the deletion.
$rngDays = $d.Content
$null = $rngDays.Find.Execute("Please enter the number of days")
$posAfterDays = $rngDays.End

$splitRange = $d.Range($posAfterDays, $posAfterDays)
$d.Bookmarks.Add("ZZTMPSPLIT", $splitRange)
$d.Bookmarks("ZZTMPSPLIT").Delete()

# Step 3: locate the boundary right after the comma, and drop the real
# "_GoBack" bookmark there (collapsed). Because "_GoBack" already exists
# elsewhere in the document, re-adding it under the same name relocates it
# here and removes it from its previous position.
$rngComma = $d.Content
$null = $rngComma.Find.Execute("Please enter the number of days,")
$posAfterComma = $rngComma.End

$bmRange = $d.Range($posAfterComma, $posAfterComma)
$d.Bookmarks.Add("_GoBack", $bmRange)
